$p = $ppt.ActivePresentation

# --- 1. Table style change -------------------------------------------------
# The deck has a single table (on the "PLENARY" slide); retarget its style
# from the custom "Table_0" style to the built-in style referenced in the
# target revision.
$oldStyleId = "{C597FC23-A9DF-4672-9FD9-CB5B9F6DD9DF}"
$newStyleId = "{5F0C2E8C-9B73-44BF-BE2B-D3CA19077E54}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($i = 1; $i -le $sl.Shapes.Count; $i++) {
        $sh = $sl.Shapes.Item($i)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Theme colour swap ---------------------------------------------------
# The deck's two embedded themes ("Office Theme" and "Integral") had their
# colour palettes swapped. The slide master/theme reachable from the object
# model is repainted with the Office Theme palette to match.
function ToRgbLong($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme 12-slot colour scheme, in COM ThemeColorScheme.Item() order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $themeColorScheme.Item($i + 1).RGB = ToRgbLong($officeThemeColors[$i])
}
